# Auto-generated script applying updated crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '26.009.83'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  +0.97%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.635.17'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  -0.95%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '215.11'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '0.504'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.993'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  -0.84%  '

$ws.Range("E8").Value = '  -0.97%  '

$ws.Range("E9").Value = '  -0.74%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '19.70'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  +0.00%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0788'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("E12").Value = '  -0.40%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '1.862.86'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +0.27%  '

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '1.586.07'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -2.89%  '

$ws.Range("E15").Value = '  -1.25%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₃0764'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  -0.06%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '63.16'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  +0.57%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '25.982.52'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("E19").Value = '  -0.88%  '

$ws.Range("E20").Value = '  -0.13%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '192.22'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  -1.15%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '10.00'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  +0.48%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '6.36'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +1.07%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '0.993'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  -0.92%  '

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '141.59'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  -0.84%  '

$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("E28").Value = '  -0.07%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '15.60'
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("E30").Value = '  +0.39%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0494'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("E32").Value = '  -0.15%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.59'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  +0.63%  '

$ws.Range("E35").Value = '  -0.21%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '0.907'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  +0.45%  '

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '1.142.85'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  +1.27%  '

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.544'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  -0.49%  '

$ws.Range("E39").Value = '  -1.82%  '

$ws.Range("E40").Value = '  +0.10%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.993'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -0.68%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '5.56'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -0.39%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '100.33'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  +0.12%  '

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '0.796'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  -1.06%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '1.772.83'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '

$ws.Range("E46").Value = '  +0.73%  '

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '55.57'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  +0.82%  '

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0517'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  +2.83%  '

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '1.45'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  +4.81%  '

$ws.Range("E50").Value = '  -0.24%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '7.60'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  +0.68%  '
